$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2567.0908
$ws.Range("J17").Value = 2567.0908
$ws.Range("L17").Value = 7701.2724
$ws.Range("N17").Value = -8037.2724

$ws.Range("H18").Value = 944.5
$ws.Range("I18").Value = 671.4
$ws.Range("K18").Value = 671.4
$ws.Range("M18").Value = -387.4

$ws.Range("H62").Value = 4622.4614
$ws.Range("I62").Value = 4700.364
$ws.Range("J62").Value = 4194
$ws.Range("K62").Value = 4700.364
$ws.Range("L62").Value = 4194
$ws.Range("M62").Value = -4076.364
$ws.Range("N62").Value = -5442

$ws.Range("H65").Value = 4622.4614
$ws.Range("I65").Value = 4700.364
$ws.Range("J65").Value = 4194
$ws.Range("K65").Value = 23501.82
$ws.Range("L65").Value = 20970
$ws.Range("M65").Value = -20381.82
$ws.Range("N65").Value = -27210

$ws.Range("H80").Value = 827.5
$ws.Range("I80").Value = 697.2222
$ws.Range("K80").Value = 2091.6666
$ws.Range("M80").Value = -1093.6666

$ws.Range("H83").Value = 827.5
$ws.Range("I83").Value = 697.2222
$ws.Range("K83").Value = 6274.999800000001
$ws.Range("M83").Value = -1282.999800000001

$ws.Range("H132").Value = 2320.279
$ws.Range("I132").Value = 1946.9524
$ws.Range("K132").Value = 5840.857199999999
$ws.Range("M132").Value = -3310.857199999999

$ws.Range("H137").Value = 6929.7896
$ws.Range("I137").Value = 8776.357
$ws.Range("K137").Value = 26329.071
$ws.Range("M137").Value = -23779.071

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4221.569
$ws.Range("I32").Value = 3051
$ws.Range("K32").Value = 3051
$ws.Range("M32").Value = -2764

$ws.Range("H45").Value = 1657.2222
$ws.Range("I45").Value = 1625.2778
$ws.Range("K45").Value = 1625.2778
$ws.Range("M45").Value = -1248.2778

$ws.Range("H61").Value = 2763.7742
$ws.Range("I61").Value = 2522.5667
$ws.Range("K61").Value = 2522.5667
$ws.Range("M61").Value = -2310.5667

$ws.Range("H122").Value = 3856.162
$ws.Range("I122").Value = 4136.3105
$ws.Range("J122").Value = 2840.625
$ws.Range("K122").Value = 12408.9315
$ws.Range("L122").Value = 8521.875
$ws.Range("M122").Value = -9958.931499999999
$ws.Range("N122").Value = -13421.875

$ws.Range("H132").Value = 21874.445
$ws.Range("I132").Value = 2406.1428
$ws.Range("J132").Value = 90013.5
$ws.Range("K132").Value = 7218.428400000001
$ws.Range("L132").Value = 270040.5
$ws.Range("M132").Value = -4688.428400000001
$ws.Range("N132").Value = -275100.5

$ws.Range("H136").Value = 2763.7742
$ws.Range("I136").Value = 2522.5667
$ws.Range("K136").Value = 7567.7001
$ws.Range("M136").Value = -5017.7001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1481.7693
$ws.Range("I134").Value = 1271.9166
$ws.Range("K134").Value = 3815.7498
$ws.Range("M134").Value = -1280.7498

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 109.55556
$ws.Range("I7").Value = 98.63636
$ws.Range("J7").Value = 126.71429
$ws.Range("K7").Value = 98.63636
$ws.Range("L7").Value = 126.71429
$ws.Range("M7").Value = 14.36364
$ws.Range("N7").Value = -352.71429

$ws.Range("H31").Value = 490549.7
$ws.Range("I31").Value = 11630.091
$ws.Range("K31").Value = 11630.091
$ws.Range("M31").Value = -11335.091

$ws.Range("H34").Value = 490549.7
$ws.Range("I34").Value = 11630.091
$ws.Range("K34").Value = 11630.091
$ws.Range("M34").Value = -11428.091

$ws.Range("H58").Value = 4189
$ws.Range("I58").Value = 4303.8335
$ws.Range("K58").Value = 4303.8335
$ws.Range("M58").Value = -4100.8335

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = ""

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

$ws.Range("H134").Value = 2755.6667
$ws.Range("I134").Value = 3167.8333
$ws.Range("K134").Value = 9503.499899999999
$ws.Range("M134").Value = -6968.499899999999

$ws.Range("H136").Value = 4189
$ws.Range("I136").Value = 4303.8335
$ws.Range("K136").Value = 12911.5005
$ws.Range("M136").Value = -10361.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 363.44446
$ws.Range("I7").Value = 275.6
$ws.Range("J7").Value = 473.25
$ws.Range("K7").Value = 826.8000000000001
$ws.Range("L7").Value = 1419.75
$ws.Range("M7").Value = -714.8000000000001
$ws.Range("N7").Value = -1643.75

$ws.Range("H11").Value = 699.75
$ws.Range("I11").Value = 433
$ws.Range("K11").Value = 1299
$ws.Range("M11").Value = -1159

$ws.Range("H33").Value = 5321.75
$ws.Range("J33").Value = 19999
$ws.Range("L33").Value = 119994
$ws.Range("N33").Value = -120560

$ws.Range("H39").Value = 3156.8333
$ws.Range("J39").Value = 4360.25
$ws.Range("L39").Value = 13080.75
$ws.Range("N39").Value = -13668.75

$ws.Range("H113").Value = 9805667
$ws.Range("I113").Value = 2728.875
$ws.Range("J113").Value = 12821955
$ws.Range("K113").Value = 8186.625
$ws.Range("L113").Value = 38465865
$ws.Range("M113").Value = -6016.625
$ws.Range("N113").Value = -38470205

$ws.Range("H122").Value = 1875
$ws.Range("I122").Value = 1875
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 16875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -14425
$ws.Range("N122").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4329.1816
$ws.Range("I80").Value = 5885.8335
$ws.Range("J80").Value = 2461.2
$ws.Range("K80").Value = 5885.8335
$ws.Range("L80").Value = 2461.2
$ws.Range("M80").Value = -4887.8335
$ws.Range("N80").Value = -4457.2

$ws.Range("H83").Value = 4329.1816
$ws.Range("I83").Value = 5885.8335
$ws.Range("J83").Value = 2461.2
$ws.Range("K83").Value = 29429.1675
$ws.Range("L83").Value = 12306
$ws.Range("M83").Value = -24437.1675
$ws.Range("N83").Value = -22290

$ws.Range("H97").Value = 1389.75
$ws.Range("I97").Value = 1586.3334
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 1586.3334
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -1090.3334
$ws.Range("N97").Value = -1792

$ws.Range("H122").Value = 1742.2941
$ws.Range("I122").Value = 1351.1875
$ws.Range("K122").Value = 4053.5625
$ws.Range("M122").Value = -1603.5625

$ws.Range("H132").Value = 4072.6
$ws.Range("J132").Value = 8564.799999999999
$ws.Range("L132").Value = 25694.4
$ws.Range("N132").Value = -30754.4

$ws.Range("H136").Value = 48800.31
$ws.Range("J136").Value = 48800.31
$ws.Range("L136").Value = 146400.93
$ws.Range("N136").Value = -151500.93

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2770.8462
$ws.Range("I16").Value = 498.08334
$ws.Range("J16").Value = 30044
$ws.Range("K16").Value = 498.08334
$ws.Range("L16").Value = 30044
$ws.Range("M16").Value = -328.08334
$ws.Range("N16").Value = -30384

$ws.Range("H40").Value = 9330.111000000001
$ws.Range("I40").Value = 7710.5713
$ws.Range("J40").Value = 14998.5
$ws.Range("K40").Value = 7710.5713
$ws.Range("L40").Value = 14998.5
$ws.Range("M40").Value = -7574.5713
$ws.Range("N40").Value = -15270.5

$ws.Range("H43").Value = 25644.334
$ws.Range("J43").Value = 19967
$ws.Range("L43").Value = 19967
$ws.Range("N43").Value = -20353

$ws.Range("H46").Value = 1632.0834
$ws.Range("I46").Value = 1129.8
$ws.Range("J46").Value = 1990.8572
$ws.Range("K46").Value = 1129.8
$ws.Range("L46").Value = 1990.8572
$ws.Range("M46").Value = -941.8
$ws.Range("N46").Value = -2366.8572

$ws.Range("H55").Value = 2167.8845
$ws.Range("I55").Value = 635.55554
$ws.Range("K55").Value = 635.55554
$ws.Range("M55").Value = -462.55554

$ws.Range("H68").Value = 4183.1333
$ws.Range("J68").Value = 5755.4443
$ws.Range("L68").Value = 5755.4443
$ws.Range("N68").Value = -7253.4443

$ws.Range("H71").Value = 4183.1333
$ws.Range("J71").Value = 5755.4443
$ws.Range("L71").Value = 28777.2215
$ws.Range("N71").Value = -36265.2215

$ws.Range("H93").Value = 5019.6665
$ws.Range("I93").Value = 2694.5715
$ws.Range("J93").Value = 13157.5
$ws.Range("K93").Value = 2694.5715
$ws.Range("L93").Value = 13157.5
$ws.Range("M93").Value = -1446.5715
$ws.Range("N93").Value = -15653.5

$ws.Range("H132").Value = 5312.579
$ws.Range("I132").Value = 4777.6816
$ws.Range("J132").Value = 7123
$ws.Range("K132").Value = 14333.0448
$ws.Range("L132").Value = 21369
$ws.Range("M132").Value = -11803.0448
$ws.Range("N132").Value = -26429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470

$ws.Range("H136").Value = 4891.2354
$ws.Range("I136").Value = 4880.7744
$ws.Range("K136").Value = 14642.3232
$ws.Range("M136").Value = -12092.3232
